# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# cecce316-...md row (row 4) and da9b6a1a-...md row (row 5) on both the
# zh-cn and de-de worksheets, reflecting the regenerated report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-20 20:18:41"
$zhcn.Range("E5").Value = "2016-03-20 20:18:41"
$zhcn.Range("H4").Value = "2016-03-20 20:19:01"
$zhcn.Range("H5").Value = "2016-03-20 20:19:01"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-20 20:18:44"
$dede.Range("E5").Value = "2016-03-20 20:18:44"
$dede.Range("H4").Value = "2016-03-20 20:19:07"
$dede.Range("H5").Value = "2016-03-20 20:19:07"
